$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 60; this shifts the existing rows 60-190
# down to 61-191, carrying their formatting (incl. the date style on
# column D) along with them.
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with this week's record.
$ws.Range("A60").Value = 4
$ws.Range("B60").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C60").Value = "Los Lagos"
$ws.Range("D60").Value = 44544
$ws.Range("D60").NumberFormat = $ws.Range("D61").NumberFormat
$ws.Range("E60").Value = 10
$ws.Range("F60").Value = 100112003
$ws.Range("G60").Value = "Ajo"
$ws.Range("H60").Value = "Chino"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 240
$ws.Range("K60").Value = 21000
$ws.Range("L60").Value = 22000
$ws.Range("M60").Value = 21500
$ws.Range("N60").Value = "$/caja 10 kilos"
$ws.Range("O60").Value = "China"
$ws.Range("P60").Value = 2150
$ws.Range("Q60").Value = 10
$ws.Range("R60").Value = "Hortaliza"
